# Rename the "Clean Training Dataset (2)" worksheet to "CSV".
# Renaming the sheet automatically updates any defined names / formulas
# that reference it (e.g. the hidden _FilterDatabase name on this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clean Training Dataset (2)")
$ws.Name = "CSV"
